# "3d case working wo scaling"
#
# Adds two more copies of the Diameter/Volume/analytical/Vf/#particles/
# Probability mini-table (the one already sitting at C28:H34) further down
# the sheet, at C37:H43 and C46:H52, each with its own tweaked inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Block 1 header (row 37)
# ---------------------------------------------------------------------------
$ws.Range("C37").Value = "Diameter"
$ws.Range("D37").Value = "Volume"
$ws.Range("E37").Value = "analytical"
$ws.Range("F37").Value = "Vf"
$ws.Range("G37").Value = "# particles"
$ws.Range("H37").Value = "Probability"

# Row 38: GC
$ws.Range("B38").Value = "GC"
$ws.Range("C38").Value = 1
$ws.Range("C38").Style = "Input"
$ws.Range("D38").Formula = "=C38^2*PI()/4"
$ws.Range("F38").Value = 0.7
$ws.Range("F38").Style = "Input"
$ws.Range("G38").Formula = "=(F38/D38)/(F`$3/D`$3)"
$ws.Range("H38").Formula = "=G38/`$G`$43"
$ws.Range("H38").NumberFormat = "0.00000"

# Row 39: fil_1 (diameter 0.42)
$ws.Range("B39").Value = "fil_1"
$ws.Range("C39").Value = 0.42
$ws.Range("C39").Style = "Input"
$ws.Range("D39").Formula = "=C39^2*PI()/4"
$ws.Range("E39").Formula = "=(C39/`$C`$30)^`$E`$1"
$ws.Range("F39").Formula = "=(E39-E40)*(1-F`$38)"
$ws.Range("G39").Formula = "=(F39/D39)/(F`$3/D`$3)"
$ws.Range("H39").Formula = "=G39/`$G`$43"
$ws.Range("H39").NumberFormat = "0.00000"

# Row 40: fil_2 (diameter 0.3)
$ws.Range("B40").Value = "fil_2"
$ws.Range("C40").Value = 0.3
$ws.Range("C40").Style = "Input"
$ws.Range("D40").Formula = "=C40^2*PI()/4"
$ws.Range("E40").Formula = "=(C40/`$C`$30)^`$E`$1"
$ws.Range("F40").Formula = "=(E40-E41)*(1-F`$39)"
$ws.Range("G40").Formula = "=(F40/D40)/(F`$3/D`$3)"
$ws.Range("H40").Formula = "=G40/`$G`$43"
$ws.Range("H40").NumberFormat = "0.00000"

# Rows 41/42: blank placeholder rows (style-only, like the template)
$ws.Range("C41").Style = "Input"
$ws.Range("H41").NumberFormat = "0.00000"

$ws.Range("C42").Style = "Input"
$ws.Range("H42").NumberFormat = "0.00000"

# Row 43: sums
$ws.Range("G34").Copy()
$ws.Range("G43").PasteSpecial(-4122)  # xlPasteFormats (font-only "Total" style)
$excel.CutCopyMode = 0
$ws.Range("G43").Formula = "=SUM(G38:G42)"
$ws.Range("H43").Formula = "=SUM(H38:H42)"
$ws.Range("H43").NumberFormat = "0.00000"

# ---------------------------------------------------------------------------
# Block 2 header (row 46)
# ---------------------------------------------------------------------------
$ws.Range("C46").Value = "Diameter"
$ws.Range("D46").Value = "Volume"
$ws.Range("E46").Value = "analytical"
$ws.Range("F46").Value = "Vf"
$ws.Range("G46").Value = "# particles"
$ws.Range("H46").Value = "Probability"

# Row 47: GC
$ws.Range("B47").Value = "GC"
$ws.Range("C47").Value = 1
$ws.Range("C47").Style = "Input"
$ws.Range("D47").Formula = "=C47^2*PI()/4"
$ws.Range("F47").Value = 0.7
$ws.Range("F47").Style = "Input"
$ws.Range("G47").Formula = "=(F47/D47)/(F`$3/D`$3)"
$ws.Range("H47").Formula = "=G47/`$G`$52"
$ws.Range("H47").NumberFormat = "0.00000"

# Row 48: fil_1 (diameter 0.3, Vf entered directly = 0.15, no analytical col)
$ws.Range("B48").Value = "fil_1"
$ws.Range("C48").Value = 0.3
$ws.Range("C48").Style = "Input"
$ws.Range("D48").Formula = "=C48^2*PI()/4"
$ws.Range("F48").Value = 0.15
$ws.Range("F48").Style = "Input"
$ws.Range("G48").Formula = "=(F48/D48)/(F`$3/D`$3)"
$ws.Range("H48").Formula = "=G48/`$G`$52"
$ws.Range("H48").NumberFormat = "0.00000"

# Row 49: fil_2 (diameter 0.2, Vf entered directly = 0.15, no analytical col)
$ws.Range("B49").Value = "fil_2"
$ws.Range("C49").Value = 0.2
$ws.Range("C49").Style = "Input"
$ws.Range("D49").Formula = "=C49^2*PI()/4"
$ws.Range("F49").Value = 0.15
$ws.Range("F49").Style = "Input"
$ws.Range("G49").Formula = "=(F49/D49)/(F`$3/D`$3)"
$ws.Range("H49").Formula = "=G49/`$G`$52"
$ws.Range("H49").NumberFormat = "0.00000"

# Rows 50/51: blank placeholder rows (style-only, like the template)
$ws.Range("C50").Style = "Input"
$ws.Range("H50").NumberFormat = "0.00000"

$ws.Range("C51").Style = "Input"
$ws.Range("H51").NumberFormat = "0.00000"

# Row 52: sums
$ws.Range("G34").Copy()
$ws.Range("G52").PasteSpecial(-4122)  # xlPasteFormats (font-only "Total" style)
$excel.CutCopyMode = 0
$ws.Range("G52").Formula = "=SUM(G47:G51)"
$ws.Range("H52").Formula = "=SUM(H47:H51)"
$ws.Range("H52").NumberFormat = "0.00000"

# --- View state: mirror the author's final selection/scroll position -------
$ws.Range("G49").Select()
$excel.ActiveWindow.ScrollRow = 19
